$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recalibrated extrapolation values after removing sub-$5 price points
# (treated as noise) from the calibration input. Columns D-H are
# recomputed for the affected maturities.

$updates = @{
    3  = @{ D = 112968.7424453835;  E = -0.03954900361912312; F = 0.1616227338645521;  G = -1.270859742402559;  H = 12.16291941517945 }
    4  = @{ D = 113787.5838315562;  E = -0.03183702583983396; F = 0.1809870267588685;  G = -1.381491378692296;  H = 11.47201725257162 }
    5  = @{ D = 114552.5222494029;  E = -0.02981518386984432; F = 0.1896711815953528;  G = -0.8650624023909192; H = 8.059864369643448 }
    6  = @{ D = 115085.5785533276;  E = -0.03692796229467925; F = 0.2240855308653421;  G = -1.420633274823637;  H = 11.41200339272869 }
    7  = @{ D = 116240.434269028;   E = -0.06546078847180872; F = 0.3562984915249692;  G = -2.108840517063721;  H = 10.1709523177209 }
    8  = @{ D = 116769.1146647389;  E = -0.04957733571440093; F = 0.2023672551671724;  G = -1.056944043934468;  H = 7.673201779996237 }
    9  = @{ D = 118345.1182495731;  E = -0.07919362700215299; F = 0.3212714182374238;  G = -1.629662202438512;  H = 10.54870345197539 }
    10 = @{ D = 119806.3499882844;  E = -0.1226619770883319;  F = 0.4447509084743156;  G = -1.885463991160664;  H = 9.429761374249825 }
    11 = @{ D = 121729.4470996111;  E = -0.1976327567257077;  F = 0.7653732441717875;  G = -2.519908501056003;  H = 11.95299193332507 }
    12 = @{ D = 112295.7497965587;  E = -0.1118196009973569;  F = 0.1618885395613342;  G = -0.8133118367470256; H = 6.827669612106231 }
    15 = @{ D = 112327.6476736164;  E = -0.08804383595053202; F = 0.1730282299836449;  G = -0.761758443683565;  H = 8.70075964703963 }
    18 = @{ D = 113160.2936710395;  E = -0.03188543207222456; F = 0.1508630638501899;  G = -0.4811946359707147; H = 6.434900174995059 }
    19 = @{ D = 113309.6656998081;  E = -0.02816224707698761; F = 0.1466049474310317;  G = -0.5929151213337743; H = 7.100451824675623 }
    20 = @{ D = 114253.6463753975;  E = -0.006348963872414823;F = 0.1402737577276532;  G = -0.4545997911603533; H = 6.186288621977192 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
}
